$wb = $excel.ActiveWorkbook

# Common formula template (as stored in row 2 of each sheet already):
#  ` "{ ""foreign"": """ & A# & """, ""grammar"": """ & B# & """, ""pronunciation"": """ & C# & """, ""meaning"": """ & D# & """ },"`
# We fill it down (as a multi-cell range assignment) so the engine creates a
# shared formula group, matching rows 3.. of each vocabulary sheet.

# --- "PAR RAPPORT À UN POINT": fill F3:F49, then update selection to F2:F49 ---
$wsPoint = $wb.Worksheets.Item("PAR RAPPORT À UN POINT")
$wsPoint.Range("F3:F49").Formula = '= "{ ""foreign"": """ & A3 & """, ""grammar"": """ & B3 & """, ""pronunciation"": """ & C3 & """, ""meaning"": """ & D3 & """ },"'
$wsPoint.Range("F2:F49").Select()
$excel.ActiveWindow.ScrollRow = 22

# --- "DISTANCE ET PROXIMITÉ": fill F3:F17, then update selection to F2:F17 ---
$wsDist = $wb.Worksheets.Item("DISTANCE ET PROXIMITÉ")
$wsDist.Range("F3:F17").Formula = '= "{ ""foreign"": """ & A3 & """, ""grammar"": """ & B3 & """, ""pronunciation"": """ & C3 & """, ""meaning"": """ & D3 & """ },"'
$wsDist.Range("F2:F17").Select()

# --- "DIRECTION ET MOUVEMENT": fill F3:F12, then update selection to F2:F12 ---
# (this sheet stays the active / tabSelected sheet, so it is selected last)
$wsDir = $wb.Worksheets.Item("DIRECTION ET MOUVEMENT")
$wsDir.Range("F3:F12").Formula = '= "{ ""foreign"": """ & A3 & """, ""grammar"": """ & B3 & """, ""pronunciation"": """ & C3 & """, ""meaning"": """ & D3 & """ },"'
$wsDir.Range("F2:F12").Select()

# --- "POUR QUALIFIER UN LIEU": only the selection changes (F2 -> F2:F7) ---
$wsQualif = $wb.Worksheets.Item("POUR QUALIFIER UN LIEU")
$wsQualif.Range("F2:F7").Select()

# --- "LOCALISATION": only the selection changes (F2 -> F2:F23) ---
$wsLoc = $wb.Worksheets.Item("LOCALISATION")
$wsLoc.Range("F2:F23").Select()

# Restore the originally active sheet / selection so "DIRECTION ET MOUVEMENT"
# stays the tab-selected sheet with F2:F12 selected.
$wsDir.Select()
$wsDir.Range("F2:F12").Select()
